{"js": "// Change the paragraph text from \"Version 2.\" to \"Version 1.\"\n// (reverting \"Wireframes version 2.\" back to \"Wireframes version 1.\").\n//\n// The edit is applied in three small, targeted text replacements instead\n// of rewriting the whole paragraph, so the existing run layout is\n// disturbed as little as possible \u2014 matching a real in-place Word edit:\n//   1. Re-type \"Version\" over the \"Versi\"/\"on\" split -> merges into one run.\n//   2. Remove the lone trailing \".\" run (now unambiguous: only one \".\"\n//      exists in the paragraph at this point).\n//   3. Re-type \" 2\" as \" 1.\" in place.\n\nconst body = context.document.body;\n\n// 1) Merge \"Versi\" + \"on\" into a single \"Version\" run.\nconst versionHits = body.search(\"Version\", { matchCase: true });\nversionHits.load(\"items\");\nawait context.sync();\nversionHits.items[0].insertText(\"Version\", \"Replace\");\nawait context.sync();\n\n// 2) Delete the separate trailing \".\" run.\nconst dotHits = body.search(\".\", { matchCase: true });\ndotHits.load(\"items\");\nawait context.sync();\ndotHits.items[0].insertText(\"\", \"Replace\");\nawait context.sync();\n\n// 3) Turn \" 2\" into \" 1.\"\nconst numHits = body.search(\" 2\", { matchCase: true });\nnumHits.load(\"items\");\nawait context.sync();\nnumHits.items[0].insertText(\" 1.\", \"Replace\");\nawait context.sync();\n", "ps1": "# Change the paragraph text from \"Version 2.\" to \"Version 1.\"\n# (reverting \"Wireframes version 2.\" back to \"Wireframes version 1.\").\n#\n# Applied as three small, targeted Find/Replace operations instead of\n# rewriting the whole paragraph, so the existing run layout is disturbed\n# as little as possible - matching a real in-place Word edit:\n#   1. Re-type \"Version\" over the \"Versi\"/\"on\" run split so it merges\n#      into a single run. (Word's Range.Text setter is a no-op when the\n#      new value equals the current text, so we first nudge it to a\n#      distinct value and then back to force the real merge.)\n#   2. Remove the lone trailing \".\" run (unambiguous at this point -\n#      only one \".\" remains in the paragraph).\n#   3. Re-type \" 2\" as \" 1.\" in place.\n\n$d = $word.ActiveDocument\n\n# 1) Merge \"Versi\" + \"on\" into a single \"Version\" run.\n$rVersion = $d.Content\n$null = $rVersion.Find.Execute(\"Version\")\n$rVersion.Text = \"Versionx\"\n\n$rVersion2 = $d.Content\n$null = $rVersion2.Find.Execute(\"Versionx\")\n$rVersion2.Text = \"Version\"\n\n# 2) Delete the separate trailing \".\" run.\n$rDot = $d.Content\n$null = $rDot.Find.Execute(\".\")\n$rDot.Text = \"\"\n\n# 3) Turn \" 2\" into \" 1.\"\n$rNum = $d.Content\n$null = $rNum.Find.Execute(\" 2\")\n$rNum.Text = \" 1.\"\n"}
